$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column E
$ws.Range("E1").Value = "Send Confirmation Email"

# Update email addresses
$ws.Range("C2").Value = "emp10@mycompany.com"
$ws.Range("C3").Value = "emp20@mycompany.com"

# Set "Send Confirmation Email" flags
$ws.Range("E2").Value = "No"
$ws.Range("E3").Value = "Yes"

# Remove hyperlink from C2, update hyperlink address/display on C3 in place
$ws.Range("C2").Hyperlinks.Delete()
$c3link = $ws.Range("C3").Hyperlinks.Item(1)
$c3link.Address = "mailto:emp20@mycompany.com"
$c3link.TextToDisplay = "emp20@mycompany.com"

# Update active cell selection
$ws.Range("E3").Select()
